# "Generate Report for Handoff"
# The localization-status report has moved past translation: the files are
# now ready to be handed off. Update the "Status" / summary cells from
# "In Translation" to "Ready for handoff" on every sheet, and refresh the
# associated timestamps to reflect the new handoff-generation time.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: zh-cn / de-de status columns (E, F) + generate-date (G) ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-28 18:56:48"

# --- zh-cn sheet: Status (C) + Latest Handoff Datetime (H) ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-28 18:56:44"

# --- de-de sheet: Status (C) + Latest Handoff Datetime (H) ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-28 18:56:48"

# --- Re-fit the Status-ish columns now that "Ready for handoff" is wider
#     than "In Translation" (mirrors the column-width bump in the report) ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
